$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()
$ws.Range("A9").Clear()

$ws.Range("B9").Value = 2035
$ws.Range("C9").Value = "2035_TM152_FBP_Plus_24_rerunTM1.5.2.5"
$ws.Range("D9").Value = "FinalBlueprint"
$ws.Range("E9").Value = "Blueprint"
$ws.Range("F9").Value = "Blueprint with TM1.5.2.5"
$ws.Range("G9").Value = "`"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION`""
$ws.Range("H9").Value = "run182"
$ws.Range("I9").Value = "current"
$ws.Range("J9").Value = "2035_TM152_FBP_Plus_24\INPUT"
$ws.Range("K9").Value = "Rerun Blueprint with TM1.5.2.5 and no changes to the network"
